# ajustes tras ejecutar check
# Updates recalculated descriptive-statistics figures (weighted means, CIs,
# sd.w, Chi2/p-value/VCramer, and group counts) across both worksheets.

$wb = $excel.ActiveWorkbook

# ---- Hoja_1 : Table3 (mpg by cyl) ----
$ws1 = $wb.Worksheets.Item("Hoja_1")

# cyl = 4 (row 5)
$ws1.Range("G5").Value = "27,9"
$ws1.Range("J5").Value = "25,7"
$ws1.Range("K5").Value = "30,1"
$ws1.Range("R5").Value = "4,7"

# cyl = 6 (row 6)
$ws1.Range("G6").Value = "18,2"
$ws1.Range("H6").Value = "-9,7"
$ws1.Range("J6").Value = "15,8"
$ws1.Range("K6").Value = "20,6"
$ws1.Range("R6").Value = "0,8"

# cyl = 8 (row 7)
$ws1.Range("G7").Value = "15,2"
$ws1.Range("H7").Value = "-12,7"
$ws1.Range("J7").Value = "14,0"
$ws1.Range("K7").Value = "16,3"
$ws1.Range("R7").Value = "1,9"

# ---- Hoja_2 : Table4 (vs by cyl) ----
$ws2 = $wb.Worksheets.Item("Hoja_2")

# vs = 0 (row 5)
$ws2.Range("D5").Value = "5,9%"
$ws2.Range("E5").Value = "9,4%"
$ws2.Range("F5").Value = 78.7
$ws2.Range("G5").Value = 1.3
$ws2.Range("H5").Value = 1.8
$ws2.Range("L5").Value = "28,50"
$ws2.Range("M5").Value = "0,0000"
$ws2.Range("N5").Value = "0,94"

# vs = 1 (row 6)
$ws2.Range("D6").Value = "94,1%"
$ws2.Range("E6").Value = "90,6%"
$ws2.Range("G6").Value = 21.2
$ws2.Range("H6").Value = 17
